# "Generate Report for Handoff"
# Regenerates the localization-status report: the in-progress status/date
# stamps move from "In Translation" to "Ready for handoff", with refreshed
# timestamps, and the Status/date columns widen to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: per-language status columns (zh-cn, de-de) and the
# "Latest HO Xliff Generate Date" timestamp.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-02 03:08:14"

# zh-cn detail sheet: Status + Latest Handoff Datetime.
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-02 03:08:09"

# de-de detail sheet: Status + Latest Handoff Datetime.
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-02 03:08:14"

# The Status/date columns grow wider to accommodate "Ready for handoff".
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
